$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44166
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 600
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 16500
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1100

# Row 3
$ws.Range("D3").Value = 44162
$ws.Range("L3").Value = "Tercera"
$ws.Range("M3").Value = 500
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 15500
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1033

# Row 4
$ws.Range("D4").Value = 44169
$ws.Range("M4").Value = 500
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 15500
$ws.Range("S4").Value = 1033

# Row 5
$ws.Range("D5").Value = 44159
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 15500
$ws.Range("P5").Value = 15750
$ws.Range("S5").Value = 1050

# Row 6
$ws.Range("D6").Value = 44187
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 350
$ws.Range("N6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("S6").Value = 1067

# Row 7
$ws.Range("D7").Value = 44187
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 13000
$ws.Range("S7").Value = 867

# Row 8
$ws.Range("D8").Value = 44194
$ws.Range("M8").Value = 300
$ws.Range("R8").Value = "Región Metropolitana"

# Row 9
$ws.Range("D9").Value = 44176
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 500
$ws.Range("N9").Value = 15000
$ws.Range("P9").Value = 15500
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 1033
